$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("reaction5")

# Clear out the old values beyond column B first
$ws.Range("C1:R1").Clear()

# Set the new values for A1 and B1
$ws.Range("A1").Value = 10
$ws.Range("B1").Value = 11
